$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing table (dates in column A, counts/sums in B:D) with
# new daily rows, from 2021-07-26 (serial 44403) through 2021-08-09
# (serial 44417), continuing the data series already present through
# row 328 (2021-07-25 / serial 44402).

# First, copy the formatting (incl. the date number-format style used in
# column A) from the last existing row down across the new rows so the
# new cells inherit the same style as the rest of the column.
$ws.Range("A328:D328").Copy($ws.Range("A329:D343"))

$dates  = 44403,44404,44405,44406,44407,44408,44409,44410,44411,44412,44413,44414,44415,44416,44417
$bVals  = 0,0,0,0,0,0,0,0,0,0,0,0,0,0,1
$cVals  = 0,0,0,0,0,0,0,0,0,0,0,0,0,0,1
$dVals  = 0,0,0,0,0,0,0,0,0,0,0,0,0,0,83.40283569641367

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 329 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
}
